# Actualizacion Formato Americano abril
# Replace the Spanish (investing.com-style) OHLC column-header labels used
# for the SP500 raw_column values with their English equivalents, matching
# the convention already used for the EURUSD / USDJPY rows further down the
# sheet. Also update the saved view state (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# raw_column (column G) header translations
$ws.Range("G2").Value  = "Date"   # was "Fecha"   (instrument=SP500, canonical_name=date)
$ws.Range("G3").Value  = "Price"  # was "Ultimo"  (instrument=SP500, canonical_name=sp500_close)
$ws.Range("G4").Value  = "Open"   # was "Apertura"(instrument=SP500, canonical_name=sp500_open)
$ws.Range("G5").Value  = "High"   # was "Maximo"  (instrument=SP500, canonical_name=sp500_high)
$ws.Range("G6").Value  = "Low"    # was "Minimo"  (instrument=SP500, canonical_name=sp500_low)
$ws.Range("G8").Value  = "Price"  # was "Ultimo"  (instrument=SP500, canonical_name=vix)
$ws.Range("G38").Value = "Price"  # was "Ultimo"  (instrument=USDJPY, canonical_name=nikkei_225_raw)

# Update the sheet's saved scroll position / selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("I39").Select()
